$wb = $excel.ActiveWorkbook

# Sheet with the raw mail log rows
$logs = $wb.Worksheets.Item("Logs")

# New row 38 data (mirrors the other log rows; column E "Antwoord" stays empty)
$logs.Range("A38").Value = "Sollicitatie marketingfunctie"
$logs.Range("B38").Value = "mailmind.test@zohomail.eu"
$logs.Range("C38").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D38").Value = "Overig"
$logs.Range("F38").Value = "2025-06-17 21:57:46"
$logs.Range("G38").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$dRange = $logs.Range("D2:D38")
$dRange.FormatConditions.Item(1).ModifyAppliesToRange($dRange)

$gRange = $logs.Range("G2:G38")
$gRange.FormatConditions.Item(1).ModifyAppliesToRange($gRange)

# Dashboard sheet: bump the "Overig" category count to include the new row
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 10
